$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ActualDate" is a new staging column inserted alphabetically among the
# existing header names. Because this sheet's header row is generated from
# an alphabetically sorted column list, adding "ActualDate" pushes every
# header from its old alphabetical position one slot to the right, and the
# header that previously occupied the last column (T) spills into a brand
# new column U.

$headers = @(
    "ActualDate",
    "ActualLabel",
    "ActualValue",
    "AgeBandBusinessKey",
    "BusinessKey",
    "CommunityTypeBusinessKey",
    "DataVersionBusinessKey",
    "DonorBusinessKey",
    "FrameworkBusinessKey",
    "GenderBusinessKey",
    "GroupBusinessKey",
    "GroupVersion",
    "IndicatorBusinessKey",
    "IndicatorValues_ID",
    "InstitutionBusinessKey",
    "LocationBusinessKey",
    "Notes",
    "OrganizationBusinessKey",
    "ReportingPeriodBusinessKey",
    "ResultAreaBusinessKey",
    "StrategicElementBusinessKey"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $headers[$i]
}
